$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 4 -> 5, Wrong penalty -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): Right total 56 -> 70, Wrong total -14 -> -16.8, Max text 42/112 -> 53.2/140
$ws.Range("B12").Value = 70
$ws.Range("C12").Value = -16.8
$ws.Range("E12").Value = "53.2/140"
